$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DeltaG")

# --- Update of databases: DeltaG0f and DeltaH0f ---
# Append a new reference row for "Sulfur (rhombic)" (S0), the standard-state
# reference phase for elemental sulfur, with a formation value of 0.
# Cells are written in the same order the new shared strings were first
# introduced (Formula, Phase, IUPAC name, Value, REF) so the shared-string
# table comes out in the same order as the authored workbook.
$ws.Range("B43").Value = "S0"
$ws.Range("C43").Value = "S"
$ws.Range("A43").Value = "Sulfur (rhombic)"
$ws.Range("D43").Value = 0
$ws.Range("E43").Value = "Kleerebezem2010"

# The data columns (B:E) use the workbook's centered number/text style,
# same as every other data row in the table.
$ws.Range("B43:E43").HorizontalAlignment = -4108   # xlCenter

# Reflect the author's final on-screen selection/scroll position.
$ws.Application.Goto($ws.Range("A18"), $false)
$ws.Range("F34").Select()
